# Apply scraped cryptos-list refresh (Thu Jul 25 06:35:45 UTC 2024).
# Rows 36/37 and 43/44 swapped rank position (coin name/link/price/volume).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '64.152.68'
$ws.Range("E2").Value = '  -2.70%  '

# Row 3
$ws.Range("D3").Value = '3.177.57'
$ws.Range("E3").Value = '  -7.78%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").Value = '''562.95'
$ws.Range("E5").Value = '  -3.62%  '

# Row 6
$ws.Range("D6").Value = '''170.74'
$ws.Range("E6").Value = '  -1.85%  '

# Row 7
$ws.Range("E7").Value = '  +0.04%  '

# Row 8
$ws.Range("D8").Value = '''0.604'
$ws.Range("E8").Value = '  +0.65%  '

# Row 9
$ws.Range("D9").Value = '3.175.16'
$ws.Range("E9").Value = '  -7.86%  '

# Row 11
$ws.Range("E11").Value = '  -4.53%  '

# Row 12
$ws.Range("D12").Value = '''0.397'
$ws.Range("E12").Value = '  -3.02%  '

# Row 13
$ws.Range("D13").Value = '3.726.45'
$ws.Range("E13").Value = '  -7.78%  '

# Row 14
$ws.Range("E14").Value = '  +1.17%  '

# Row 15
$ws.Range("D15").Value = '''27.43'
$ws.Range("E15").Value = '  -5.01%  '

# Row 16
$ws.Range("D16").Value = '64.163.33'
$ws.Range("E16").Value = '  -2.67%  '

# Row 17
$ws.Range("E17").Value = '  -4.97%  '

# Row 18
$ws.Range("D18").Value = '3.176.86'
$ws.Range("E18").Value = '  -7.71%  '

# Row 19
$ws.Range("D19").Value = '''5.70'
$ws.Range("E19").Value = '  -4.07%  '

# Row 20
$ws.Range("D20").Value = '''13.10'
$ws.Range("E20").Value = '  -4.92%  '

# Row 21
$ws.Range("D21").Value = '''352.97'
$ws.Range("E21").Value = '  -4.22%  '

# Row 22
$ws.Range("D22").Value = '''7.21'
$ws.Range("E22").Value = '  -4.93%  '

# Row 23
$ws.Range("D23").Value = '''1.00'

# Row 24
$ws.Range("D24").Value = '''68.82'
$ws.Range("E24").Value = '  -4.71%  '

# Row 25
$ws.Range("D25").Value = '''0.505'
$ws.Range("E25").Value = '  -4.57%  '

# Row 26
$ws.Range("E26").Value = '  -3.06%  '

# Row 27
$ws.Range("D27").Value = '''9.68'
$ws.Range("E27").Value = '  -0.33%  '

# Row 28
$ws.Range("D28").Value = '''0.174'
$ws.Range("E28").Value = '  -2.24%  '

# Row 29
$ws.Range("E29").Value = '  -0.03%  '

# Row 30
$ws.Range("D30").Value = '''5.69'
$ws.Range("E30").Value = '  -1.27%  '

# Row 31
$ws.Range("D31").Value = '''1.00'
$ws.Range("E31").Value = '  +0.07%  '

# Row 32
$ws.Range("D32").Value = '''1.90'
$ws.Range("E32").Value = '  -3.97%  '

# Row 33
$ws.Range("D33").Value = '''22.12'
$ws.Range("E33").Value = '  -6.25%  '

# Row 34
$ws.Range("D34").Value = '''6.66'
$ws.Range("E34").Value = '  -4.49%  '

# Row 35
$ws.Range("E35").Value = '  -5.52%  '

# Row 36
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '''1.44'
$ws.Range("E36").Value = '  -5.66%  '

# Row 37
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").Value = '''155.31'
$ws.Range("E37").Value = '  -3.51%  '

# Row 38
$ws.Range("D38").Value = '''0.807'
$ws.Range("E38").Value = '  -8.40%  '

# Row 39
$ws.Range("D39").Value = '''25.91'
$ws.Range("E39").Value = '  -9.08%  '

# Row 40
$ws.Range("E40").Value = '  -5.00%  '

# Row 41
$ws.Range("D41").Value = '''2.50'
$ws.Range("E41").Value = '  -5.21%  '

# Row 42
$ws.Range("D42").Value = '2.652.43'
$ws.Range("E42").Value = '  -4.08%  '

# Row 43
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = '''4.17'
$ws.Range("E43").Value = '  -6.11%  '

# Row 44
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '''5.99'
$ws.Range("E44").Value = '  -7.47%  '

# Row 45
$ws.Range("D45").Value = '''330.39'
$ws.Range("E45").Value = '  +1.70%  '

# Row 46
$ws.Range("D46").Value = '''0.0655'
$ws.Range("E46").Value = '  -3.62%  '

# Row 47
$ws.Range("E47").Value = '  -4.07%  '

# Row 48
$ws.Range("D48").Value = '''23.73'
$ws.Range("E48").Value = '  -3.26%  '

# Row 49
$ws.Range("D49").Value = '''0.0272'
$ws.Range("E49").Value = '  -6.15%  '

# Row 50
$ws.Range("E50").Value = '  -1.51%  '

# Row 51
$ws.Range("E51").Value = '  -0.10%  '
